$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '68.112.38'
$ws.Range('E2').Value = '  -1.00%  '
Set-TextCell $ws.Range('D3') '2.642.40'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextCell $ws.Range('D5') '596.51'
$ws.Range('E5').Value = '  -0.67%  '
Set-TextCell $ws.Range('D6') '155.62'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.10%  '
Set-TextCell $ws.Range('D9') '0.141'
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  -0.36%  '
Set-TextCell $ws.Range('D13') '27.94'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('E14').Value = '  -0.37%  '
Set-TextCell $ws.Range('D15') '3.124.52'
$ws.Range('E15').Value = '  -0.25%  '
Set-TextCell $ws.Range('D16') '68.145.62'
$ws.Range('E16').Value = '  -0.83%  '
Set-TextCell $ws.Range('D17') '2.635.40'
$ws.Range('E17').Value = '  -0.17%  '
Set-TextCell $ws.Range('D18') '11.34'
$ws.Range('E18').Value = '  -0.45%  '
Set-TextCell $ws.Range('D19') '362.98'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('E20').Value = '  -1.21%  '
Set-TextCell $ws.Range('D21') '4.40'
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('E22').Value = '  -2.78%  '
Set-TextCell $ws.Range('D23') '2.06'
Set-TextCell $ws.Range('D24') '74.73'
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  -3.56%  '
Set-TextCell $ws.Range('D27') '2.773.80'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -2.31%  '
Set-TextCell $ws.Range('D29') '1.00'
$ws.Range('E29').Value = '  -0.13%  '
Set-TextCell $ws.Range('D30') '554.57'
$ws.Range('E30').Value = '  -4.93%  '
Set-TextCell $ws.Range('D31') '8.00'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('E34').Value = '  -2.33%  '
Set-TextCell $ws.Range('D35') '1.00'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -0.82%  '
Set-TextCell $ws.Range('D37') '161.12'
$ws.Range('E37').Value = '  +0.40%  '
Set-TextCell $ws.Range('D38') '19.44'
$ws.Range('E38').Value = '  +0.55%  '
Set-TextCell $ws.Range('D39') '0.371'
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('E40').Value = '  -3.61%  '
Set-TextCell $ws.Range('D41') '5.31'
$ws.Range('E41').Value = '  -1.46%  '
Set-TextCell $ws.Range('D42') '0.0₆0338'
$ws.Range('E42').Value = '  +4.77%  '
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('E44').Value = '  -2.28%  '
Set-TextCell $ws.Range('B46') 'Aave'
Set-TextCell $ws.Range('C46') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws.Range('D46') '159.32'
$ws.Range('E46').Value = '  +1.96%  '
Set-TextCell $ws.Range('B47') 'Filecoin'
Set-TextCell $ws.Range('C47') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range('D47') '3.72'
$ws.Range('E47').Value = '  -0.45%  '
Set-TextCell $ws.Range('B48') 'InjectiveProtocol'
Set-TextCell $ws.Range('C48') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range('D48') '21.98'
$ws.Range('E48').Value = '  -0.41%  '
Set-TextCell $ws.Range('B49') 'Optimism'
Set-TextCell $ws.Range('C49') 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
Set-TextCell $ws.Range('D49') '1.69'
$ws.Range('E49').Value = '  -1.44%  '
Set-TextCell $ws.Range('B50') 'Cronos'
Set-TextCell $ws.Range('C50') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range('D50') '0.0782'
$ws.Range('E50').Value = '  -0.01%  '
Set-TextCell $ws.Range('B51') 'Mantle'
Set-TextCell $ws.Range('C51') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws.Range('D51') '0.614'
$ws.Range('E51').Value = '  -0.63%  '
